$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "le poids du département dans sa région et et le poids",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "le poids du département dans sa région et le poids", 2
)
